$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the scheduled GitHub Actions refresh.
# Plain text / non-numeric-looking values are assigned directly.
# Values that look like plain numbers are assigned with a leading apostrophe
# (forces Excel to keep them as literal text, matching the source data which
# stores prices as strings) and then the cell style is reset back to Normal so
# no stray numeric formatting is left behind.

$ws.Range("D2").Value = "28.026.05"
$ws.Range("E2").Value = "  +4.95%  "
$ws.Range("D3").Value = "1.781.70"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'244.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4923"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2676"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").Value = "'0.06266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "1.780.36"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("D11").Value = "'16.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").Value = "'0.07033"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'0.6281"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "'4.641"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'80.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "'0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.979.24"
$ws.Range("E17").Value = "  +5.56%  "
$ws.Range("D18").Value = "'0.9992"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'0.000007223"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +4.79%  "
$ws.Range("D21").Value = "2.013.73"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").Value = "'4.575"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "'8.706"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "'5.231"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "'142.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "'15.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'1.862"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.05%  "
$ws.Range("D28").Value = "'109.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("D29").Value = "'1.388"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'4.198"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.05%  "
$ws.Range("D31").Value = "'0.08284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "'3.796"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "'0.04891"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.66%  "
$ws.Range("D34").Value = "'1.075"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.17%  "
$ws.Range("D35").Value = "'2.612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'0.6529"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D37").Value = "'0.9482"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "'2.583"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.97%  "
$ws.Range("D39").Value = "'2.054"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").Value = "'5.994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.06%  "
$ws.Range("D41").Value = "'0.01554"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'100.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'0.3993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").Value = "'7.191"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("D46").Value = "'0.1201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").Value = "'0.05423"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "'8.024"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "'1.304"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.52%  "
$ws.Range("D50").Value = "'30.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'52.94"
$ws.Range("D51").Style = "Normal"
